# Manual testing for submission - sequence diagram tweaks on the "undo" slide
# (slide 10: the "execute(undo)" / "retreivePreviousStateInUndoStack" diagram).
#
# Shape.Left/.Top/.Width/.Height are COM `Single` (32-bit float) properties,
# so EMU values are recovered as floor(Single(points) * 12700). A tiny
# sub-EMU epsilon is added before the EMU->points conversion so the
# round-trip lands back on the exact target EMU instead of one under it.
$EmuPerPoint = 12700.0
$Epsilon = 0.75 / $EmuPerPoint
function EmuToPt($emu) { return ($emu / $EmuPerPoint) + $Epsilon }

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

# --- Rectangle 62 (the ":Model" box atop its lifeline) ---
$rect62 = $s.Shapes.Item(8)
$rect62.Left = EmuToPt 7239000
$rect62.Top  = EmuToPt 579075

# --- Straight Connector 13 (the ":Model" lifeline) ---
$conn13 = $s.Shapes.Item(9)
$conn13.Left = EmuToPt 7785817
$conn13.Top  = EmuToPt 942746

# --- Rectangle 14 (the activation bar on the ":Model" lifeline) ---
$rect14 = $s.Shapes.Item(10)
$rect14.Left = EmuToPt 7713809
$rect14.Top  = EmuToPt 1509557

# --- TextBox 18: merge the split "execute" / "(“undo”)" runs into a single run ---
$execBox = $s.Shapes.Item(14)
$execRange = $execBox.TextFrame.TextRange
$whole = $execRange.Characters(1, $execRange.Length)
$whole.Text = "execute(“undo”)"

# --- TextBox 20: reposition/resize the "retreivePreviousStateInUndoStack" label ---
$undoStackBox = $s.Shapes.Item(16)
$undoStackBox.Left   = EmuToPt 4642220
$undoStackBox.Top    = EmuToPt 1847870
$undoStackBox.Width  = EmuToPt 2999581
$undoStackBox.Height = EmuToPt 215443
